# Apply updated "dSF" (column F) values for the rows that were repulled/recomputed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    4  = -8
    6  = -4
    8  = -9
    10 = -7
    11 = -6
    12 = -4
    13 = -4
    14 = -11
    17 = 0
    26 = 1
    33 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
